$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 125808.625
$ws.Range("I28").Value = 154495.23
$ws.Range("J28").Value = 1500
$ws.Range("K28").Value = 154495.23
$ws.Range("L28").Value = 1500
$ws.Range("M28").Value = -154010.23
$ws.Range("N28").Value = -2470

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 20913.6
$ws.Range("I38").Value = 1142
$ws.Range("J38").Value = 100000
$ws.Range("K38").Value = 3426
$ws.Range("L38").Value = 300000
$ws.Range("M38").Value = -3054
$ws.Range("N38").Value = -300744

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3575.742
$ws.Range("I40").Value = 1858.3334
$ws.Range("J40").Value = 3987.92
$ws.Range("K40").Value = 1858.3334
$ws.Range("L40").Value = 3987.92
$ws.Range("M40").Value = -1683.3334
$ws.Range("N40").Value = -4337.92

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 6900.273
$ws.Range("I64").Value = 4650.3335
$ws.Range("J64").Value = 9600.200000000001
$ws.Range("K64").Value = 4650.3335
$ws.Range("L64").Value = 9600.200000000001
$ws.Range("M64").Value = -4402.3335
$ws.Range("N64").Value = -10096.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 6900.273
$ws.Range("I67").Value = 4650.3335
$ws.Range("J67").Value = 9600.200000000001
$ws.Range("K67").Value = 4650.3335
$ws.Range("L67").Value = 9600.200000000001
$ws.Range("M67").Value = -3792.3335
$ws.Range("N67").Value = -11316.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 14706591
$ws.Range("I92").Value = 17857882
$ws.Range("K92").Value = 17857882
$ws.Range("M92").Value = -17856634

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 32169.143
$ws.Range("I137").Value = 45628.79
$ws.Range("K137").Value = 136886.37
$ws.Range("M137").Value = -134336.37

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3145.2307
$ws.Range("I2").Value = 2932.3333
$ws.Range("J2").Value = 3624.25
$ws.Range("K2").Value = 2932.3333
$ws.Range("L2").Value = 3624.25
$ws.Range("M2").Value = -2819.3333
$ws.Range("N2").Value = -3850.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2315.4167
$ws.Range("I45").Value = 1285.625
$ws.Range("K45").Value = 1285.625
$ws.Range("M45").Value = -908.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 3145.2307
$ws.Range("I116").Value = 2932.3333
$ws.Range("J116").Value = 3624.25
$ws.Range("K116").Value = 2932.3333
$ws.Range("L116").Value = 3624.25
$ws.Range("M116").Value = -638.3332999999998
$ws.Range("N116").Value = -8212.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 403451.56
$ws.Range("I132").Value = 558626.5
$ws.Range("K132").Value = 1675879.5
$ws.Range("M132").Value = -1673349.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3145.2307
$ws.Range("I3").Value = 2932.3333
$ws.Range("J3").Value = 3624.25
$ws.Range("K3").Value = 2932.3333
$ws.Range("L3").Value = 3624.25
$ws.Range("M3").Value = -2818.3333
$ws.Range("N3").Value = -3852.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H38").Value = 98830
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 98830
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 98830
$ws.Range("M38").ClearContents()
$ws.Range("N38").Value = -99662

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1702829.4
$ws.Range("I134").Value = 1985275.9
$ws.Range("K134").Value = 5955827.699999999
$ws.Range("M134").Value = -5953292.699999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4652.5537
$ws.Range("I31").Value = 2164.1516
$ws.Range("J31").Value = 7218.7188
$ws.Range("K31").Value = 2164.1516
$ws.Range("L31").Value = 7218.7188
$ws.Range("M31").Value = -1869.1516
$ws.Range("N31").Value = -7808.7188

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4652.5537
$ws.Range("I34").Value = 2164.1516
$ws.Range("J34").Value = 7218.7188
$ws.Range("K34").Value = 2164.1516
$ws.Range("L34").Value = 7218.7188
$ws.Range("M34").Value = -1962.1516
$ws.Range("N34").Value = -7622.7188

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2945515.8
$ws.Range("I122").Value = 4549373
$ws.Range("J122").Value = 5110.5
$ws.Range("K122").Value = 13648119
$ws.Range("L122").Value = 15331.5
$ws.Range("M122").Value = -13645669
$ws.Range("N122").Value = -20231.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3584.8108
$ws.Range("I132").Value = 3379.5715
$ws.Range("K132").Value = 10138.7145
$ws.Range("M132").Value = -7608.7145

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3735
$ws.Range("I134").Value = 4203
$ws.Range("J134").Value = 3111
$ws.Range("K134").Value = 12609
$ws.Range("L134").Value = 9333
$ws.Range("M134").Value = -10074
$ws.Range("N134").Value = -14403

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 20408.809
$ws.Range("J70").Value = 5956.533
$ws.Range("L70").Value = 5956.533
$ws.Range("N70").Value = -6496.533

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 20408.809
$ws.Range("J73").Value = 5956.533
$ws.Range("L73").Value = 5956.533
$ws.Range("N73").Value = -7828.533

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1627.2354
$ws.Range("I102").Value = 1649.5927
$ws.Range("K102").Value = 1649.5927
$ws.Range("M102").Value = -27.59269999999992

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 792.0909
$ws.Range("I122").Value = 792.0909
$ws.Range("K122").Value = 2376.2727
$ws.Range("M122").Value = 73.72730000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 37040000
$ws.Range("I40").Value = 41668124
$ws.Range("K40").Value = 41668124
$ws.Range("M40").Value = -41667988

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2500
$ws.Range("J93").Value = 3800
$ws.Range("L93").Value = 3800
$ws.Range("N93").Value = -6296

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 7856.16
$ws.Range("I122").Value = 7919.4736
$ws.Range("J122").Value = 7655.6665
$ws.Range("K122").Value = 23758.4208
$ws.Range("L122").Value = 22966.9995
$ws.Range("M122").Value = -21308.4208
$ws.Range("N122").Value = -27866.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 271029.9
$ws.Range("I132").Value = 299573.7
$ws.Range("J132").Value = 6999.75
$ws.Range("K132").Value = 898721.1000000001
$ws.Range("L132").Value = 20999.25
$ws.Range("M132").Value = -896191.1000000001
$ws.Range("N132").Value = -26059.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1734.909
$ws.Range("I136").Value = 1826.762
$ws.Range("K136").Value = 5480.286
$ws.Range("M136").Value = -2930.286

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4306.1577
$ws.Range("I126").Value = 5343.1665
$ws.Range("J126").Value = 2528.4285
$ws.Range("K126").Value = 16029.4995
$ws.Range("L126").Value = 7585.2855
$ws.Range("M126").Value = -13559.4995
$ws.Range("N126").Value = -12525.2855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 23997.262
$ws.Range("I132").Value = 27404.334
$ws.Range("K132").Value = 82213.00199999999
$ws.Range("M132").Value = -79683.00199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 38845.89
$ws.Range("I136").Value = 1108.762
$ws.Range("K136").Value = 3326.286
$ws.Range("M136").Value = -776.2860000000001

